$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# D-column values are forced to Text so numeric-looking strings (e.g. "1.000")
# keep their exact original formatting instead of being coerced to numbers.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.640.48'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E2").Value = '  +4.21%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.746.17'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E3").Value = '  +4.53%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9996'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '247.20'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +3.28%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.02%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4804'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.20%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2705'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E8").Value = '  +2.71%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06250'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.15%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.746.13'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +4.49%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07109'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.57%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '15.86'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +6.82%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.6189'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +4.67%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '4.507'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.65%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '77.28'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = '  +2.67%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '1.0000'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.01%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '26.635.81'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E17").Value = '  +4.21%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.10%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000006904'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.04%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '11.72'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.99%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '1.970.35'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +4.57%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.638'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +4.44%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '8.868'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.36%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '5.354'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.37%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '136.76'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  +2.88%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '1.826'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +5.85%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.419'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = '  +2.06%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '107.79'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.91%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '4.028'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("E31").Value = '  +3.12%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.07900'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.05%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.04578'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +8.39%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.615'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.17%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.6371'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value = '  +4.39%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.9988'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +4.80%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.9535'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +11.24%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '114.01'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = '  +18.19%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.484'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -4.29%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.973'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +5.17%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.29%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.01519'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +2.98%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '5.718'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +17.49%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.3920'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +3.96%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '6.722'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = '  +8.26%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.1204'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +7.79%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.05328'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.30%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.957'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E48").Value = '  +7.84%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '30.84'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.24%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.3454'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +3.44%  '
$ws.Range("E51").Value = '  +3.36%  '
